$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): copy the existing header style (bold, centered,
# bordered) from L1 onto the new header cells, then set their text.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$ws.Range("M1").Value() = "renewd"
$ws.Range("N1").Value() = "PlanID"
$ws.Range("O1").Value() = "iteration"

# New data columns for every existing data row (rows 2-21)
for ($r = 2; $r -le 21; $r++) {
  $ws.Cells.Item($r, 13).Value() = "after"
  $ws.Cells.Item($r, 14).Value() = 20120398
  $ws.Cells.Item($r, 15).Value() = 1
}
